$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.02"
$ws.Range("E2").Value = "'-0.96%"
$ws.Range("D3").Value = "'37.75"
$ws.Range("E3").Value = "'0.18%"
$ws.Range("D4").Value = "'5.115"
$ws.Range("E4").Value = "'0.55%"
$ws.Range("D5").Value = "'0.07848"
$ws.Range("E5").Value = "'0.96%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.900"
$ws.Range("E6").Value = "'0.06%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.234"
$ws.Range("E7").Value = "'0.32%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.991"
$ws.Range("E8").Value = "'2.00%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9338"
$ws.Range("E9").Value = "'1.93%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1075"
$ws.Range("E10").Value = "'-10.43%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1930"
$ws.Range("E11").Value = "'0.78%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09028"
$ws.Range("E12").Value = "'-2.62%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03338"
$ws.Range("E13").Value = "'-1.99%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09586"
$ws.Range("E14").Value = "'-1.08%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001384"
$ws.Range("E15").Value = "'1.52%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005804"
$ws.Range("E16").Value = "'0.07%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.615"
$ws.Range("E17").Value = "'1.63%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.428"
$ws.Range("E18").Value = "'1.86%"
$ws.Range("E19").Value = "'1.13%"
$ws.Range("D20").Value = "'6.212"
$ws.Range("E20").Value = "'23.39%"
$ws.Range("E21").Value = "'1.04%"
$ws.Range("E22").Value = "'-10.49%"
$ws.Range("D23").Value = "'0.04401"
$ws.Range("E23").Value = "'0.71%"
$ws.Range("D24").Value = "'0.001232"
$ws.Range("E24").Value = "'1.58%"
$ws.Range("D25").Value = "'0.004565"
$ws.Range("E25").Value = "'7.27%"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'0.17%"
$ws.Range("D39").Value = "'0.02219"
$ws.Range("E39").Value = "'4.84%"
$ws.Range("D40").Value = "'0.05057"
$ws.Range("E40").Value = "'2.17%"
$ws.Range("D41").Value = "'0.007472"
$ws.Range("E41").Value = "'-2.24%"
$ws.Range("D42").Value = "'0.1349"
$ws.Range("E42").Value = "'0.49%"
$ws.Range("D43").Value = "'0.008740"
$ws.Range("E43").Value = "'-11.93%"
$ws.Range("D44").Value = "'0.002113"
$ws.Range("E44").Value = "'2.61%"
$ws.Range("D45").Value = "'0.007972"
$ws.Range("E45").Value = "'-9.36%"
$ws.Range("D46").Value = "'0.00006531"
$ws.Range("E46").Value = "'-2.03%"
$ws.Range("D48").Value = "'0.002861"
$ws.Range("E48").Value = "'-5.88%"
$ws.Range("E49").Value = "'-40.74%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D51").Value = "'0.0002001"
